# LDLC price-tracking workbook: append a new scrape snapshot column.
#
# Layout (before edit):
#   A            reference
#   B .. AA      one column per timestamped price scrape (27 snapshot cols)
#   AB           nom (product name)
#   AC           url_produit (product URL)
#
# This edit inserts a brand-new snapshot column right after the last
# existing snapshot column (AA), labelled with the new scrape timestamp.
# Its per-row values are carried forward from the previous (AA) snapshot
# column. "nom" and "url_produit" shift one column to the right (AB->AC,
# AC->AD) as a natural consequence of the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 205
$oldLastSnapshotCol = 27   # column AA (27th column) - last price/timestamp column before the edit
$newSnapshotCol = 28       # column AB (28th column) - the newly inserted column

# Shift AB:AD -> AC:AD (nom, url_produit) one column to the right and
# open up a fresh, blank column AB for the new snapshot.
$ws.Columns("AB:AB").Insert()

# Header: new snapshot column gets the new scrape timestamp label.
$ws.Cells.Item(1, $newSnapshotCol).Value = "2026-01-28 21:19:29"

# Data rows: copy each row's previous (AA) snapshot value into the new
# (AB) column so the new snapshot starts from the last known price.
for ($r = 2; $r -le $lastRow; $r++) {
    $srcCell = $ws.Cells.Item($r, $oldLastSnapshotCol)
    $dstCell = $ws.Cells.Item($r, $newSnapshotCol)
    $val = $srcCell.Value2
    if ($val -ne $null -and $val -ne "") {
        $dstCell.Value = $val
    }
}
